$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D7"   = -7.067
    "A8"   = -22.188
    "A10"  = -21.856
    "A12"  = -21.58900000000001
    "D14"  = -7.776999999999999
    "D15"  = -8.196000000000002
    "A18"  = -22.269
    "D18"  = -8.366
    "D20"  = -7.37
    "A25"  = -21.691
    "D29"  = -7.292
    "D30"  = -7.178999999999999
    "D31"  = -8.023
    "D35"  = -7.582000000000001
    "A37"  = -20.105
    "D40"  = -7.582000000000001
    "D44"  = -7.421000000000001
    "D50"  = -8.104999999999999
    "D54"  = -8.016000000000002
    "A55"  = -22.279
    "A68"  = -21.534
    "D68"  = -6.778
    "D76"  = -7.672000000000001
    "A77"  = -20.677
    "A78"  = -20.22
    "A79"  = -21.791
    "A80"  = -20.272
    "A81"  = -21.797
    "A82"  = -22.261
    "A84"  = -22.143
    "D87"  = -8.297000000000001
    "D88"  = -8.259
    "D92"  = -6.584000000000001
    "D96"  = -7.267
    "D98"  = -8.404
    "A101" = -21.183
    "D101" = -7.632
    "A102" = -20.418
    "D102" = -8.090999999999999
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
